$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.918.53"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.042.60"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.87%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.20"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.51"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +7.00%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.98"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0786"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.32%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "16.16"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.341.65"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.806"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.61"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +6.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.046.76"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.857.85"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.84"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +16.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "74.85"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0902"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.94%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.81"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.49%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.38"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.48%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +12.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.97"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.28"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.92%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.43%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.73"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0620"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.47"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.59%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0863"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.22"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.11%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +13.99%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.83"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.13"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.98"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.51%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.65"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +15.14%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.283.17"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.88"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.76"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.228.24"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.10%  "
